{"js": "// Delete the bullet paragraph \"Association rules is an unsupervised learning method.\"\n// from the \"Data Pre-Processing Techniques\" list.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst target = \"Association rules is an unsupervised learning method.\";\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.trim() === target) {\n    para.delete();\n  }\n}\nawait context.sync();\n", "ps1": "$doc = $word.ActiveDocument\nforeach ($para in $doc.Paragraphs) {\n    if ($para.Range.Text.TrimEnd(\"`r`a\").Trim() -eq \"Association rules is an unsupervised learning method.\") {\n        $para.Range.Delete()\n        break\n    }\n}\n"}
